$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for rows 2-6, columns B:H (A=index, I/J unchanged)
$data = @{
    2 = @{ B=7; C=1; D=4; E=4; F=-3; G=3; H=34 }
    3 = @{ B=6; C=3; D=5; E=8; F=-1; G=5; H=56 }
    4 = @{ B=8; C=4; D=4; E=6; F=-4; G=2; H=23 }
    5 = @{ B=5; C=2; D=3; E=6; F=-2; G=4; H=45 }
    6 = @{ B=7; C=2; D=2; E=3; F=-5; G=1; H=12 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
